# "Neue Messung mit Standardabweichung"
#
# - Updates the timestamp in A1 (both sheets)
# - Renames "Mittelwert" -> "Arth. Mittelwert" (both sheets, column F header)
# - Adds a new "Standardabweichung" column (H) after the mean column (both sheets)
# - Refreshes the measured figures with a new measurement run
# - Adds a "Messreihen" row label (A9) on the "Performanz Messung" sheet
# - Re-freezes the header on "Warmlaufen" to cover the first two rows

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Warmlaufen")
$ws2 = $wb.Worksheets.Item("Performanz Messung")

# ---------------------------------------------------------------------------
# Shared header text updates (A1 timestamp, F2 "Mittelwert" -> "Arth. Mittelwert")
# ---------------------------------------------------------------------------
foreach ($ws in @($ws1, $ws2)) {
    $ws.Range("A1").Value = "02.08.2014 um 22:16 Uhr"
    $ws.Range("F2").Value = "Arth. Mittelwert"

    # New "Standardabweichung" column header, styled like the neighbouring
    # header cells (bold, centered) by copying the format from G2.
    $ws.Range("G2").Copy()
    $ws.Range("H2").PasteSpecial(-4122)
    $ws.Range("H2").Value = "Standardabweichung"
}

# ---------------------------------------------------------------------------
# "Warmlaufen" sheet: new measurement run, plus stddev column
# ---------------------------------------------------------------------------

# Give the new G column (std. dev) the same number format as column F.
$ws1.Range("F3").Copy()
$ws1.Range("G3:G6").PasteSpecial(-4122)

# Row 3 - ByHand
$ws1.Range("H3").Value = 2.0
$ws1.Range("G3").Value = 0.483046

# Row 4 - Dozer
$ws1.Range("C4").Value = 950.0
$ws1.Range("D4").Value = 28.0
$ws1.Range("F4").Value = 95.0
$ws1.Range("G4").Value = 119.139
$ws1.Range("H4").Value = 28.0

# Row 5 - Orika
$ws1.Range("C5").Value = 289.0
$ws1.Range("D5").Value = 4.0
$ws1.Range("F5").Value = 28.9
$ws1.Range("G5").Value = 74.8813
$ws1.Range("H5").Value = 5.0

# Row 6 - MapStruct
$ws1.Range("C6").Value = 22.0
$ws1.Range("F6").Value = 2.2
$ws1.Range("G6").Value = 0.788811
$ws1.Range("H6").Value = 2.0

# ---------------------------------------------------------------------------
# "Performanz Messung" sheet: new measurement run, plus stddev column
# ---------------------------------------------------------------------------

$ws2.Range("F3").Copy()
$ws2.Range("G3:G6").PasteSpecial(-4122)

# Row 3 - ByHand
$ws2.Range("B3").Value = 51774.0
$ws2.Range("C3").Value = 980.0
$ws2.Range("E3").Value = 6.0
$ws2.Range("F3").Value = 0.0189284
$ws2.Range("G3").Value = 0.147176
$ws2.Range("H3").Value = 0.0

# Row 4 - Dozer
$ws2.Range("B4").Value = 51774.0
$ws2.Range("C4").Value = 532851.0
$ws2.Range("D4").Value = 9.0
$ws2.Range("E4").Value = 34.0
$ws2.Range("F4").Value = 10.2919
$ws2.Range("G4").Value = 1.1458
$ws2.Range("H4").Value = 10.0

# Row 5 - Orika
$ws2.Range("B5").Value = 51774.0
$ws2.Range("C5").Value = 27363.0
$ws2.Range("F5").Value = 0.528509
$ws2.Range("G5").Value = 0.547114
$ws2.Range("H5").Value = 1.0

# Row 6 - MapStruct
$ws2.Range("B6").Value = 51774.0
$ws2.Range("C6").Value = 14732.0
$ws2.Range("E6").Value = 4.0
$ws2.Range("F6").Value = 0.284544
$ws2.Range("G6").Value = 0.458884
$ws2.Range("H6").Value = 0.0

# New row label "Messreihen" in A9, styled like the other bold header cells.
$ws2.Range("A2").Copy()
$ws2.Range("A9").PasteSpecial(-4122)
$ws2.Range("A9").Value = "Messreihen"

# Row 10 - ByHand
$ws2.Range("C10").Value = 0.59
$ws2.Range("D10").Value = 0.068
$ws2.Range("E10").Value = 0.0234
$ws2.Range("F10").Value = 0.0213
$ws2.Range("G10").Value = 0.0200667
$ws2.Range("H10").Value = 0.01955
$ws2.Range("I10").Value = 0.01906

# Row 11 - Dozer
$ws2.Range("B11").Value = 23.4
$ws2.Range("C11").Value = 18.94
$ws2.Range("D11").Value = 11.18
$ws2.Range("E11").Value = 10.3857
$ws2.Range("F11").Value = 10.3307
$ws2.Range("G11").Value = 10.3024
$ws2.Range("H11").Value = 10.2938
$ws2.Range("I11").Value = 10.2909

# Row 12 - Orika
$ws2.Range("B12").Value = 3.9
$ws2.Range("C12").Value = 2.78
$ws2.Range("D12").Value = 0.754
$ws2.Range("E12").Value = 0.5519
$ws2.Range("F12").Value = 0.54005
$ws2.Range("G12").Value = 0.531933
$ws2.Range("H12").Value = 0.529425
$ws2.Range("I12").Value = 0.52892

# Row 13 - MapStruct
$ws2.Range("C13").Value = 0.99
$ws2.Range("D13").Value = 0.391
$ws2.Range("E13").Value = 0.3009
$ws2.Range("F13").Value = 0.28565
$ws2.Range("G13").Value = 0.286933
$ws2.Range("H13").Value = 0.286575
$ws2.Range("I13").Value = 0.2844

# ---------------------------------------------------------------------------
# Freeze the first two rows on "Warmlaufen" (used to be only the first row)
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("B3").Select()
$excel.ActiveWindow.FreezePanes = $true
